# Auto-generated: apply meteocat daily-summary refresh (data + timestamps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = "2026-02-15 20:18:31"
$ws.Range("H2").Value2 = "'70%"
$ws.Range("I2").Value2 = "2.0 mm"
$ws.Range("K2").Value2 = "4.8 MJ/m2"
$ws.Range("M2").Value2 = "4.2 °C 16:26 TU"
$ws.Range("O2").Value2 = "0.4 °C"
$ws.Range("E3").Value2 = "2026-02-15 20:18:33"
$ws.Range("I3").Value2 = "1.8 mm"
$ws.Range("O3").Value2 = "-5.1 °C"
$ws.Range("E4").Value2 = "2026-02-15 20:18:36"
$ws.Range("H4").Value2 = "'71%"
$ws.Range("E5").Value2 = "2026-02-15 20:18:43"
$ws.Range("I5").Value2 = "5.9 mm"
$ws.Range("L5").Value2 = "40.0 km/h - 312º 19:59 TU"
$ws.Range("E6").Value2 = "2026-02-15 20:18:45"
$ws.Range("H6").Value2 = "'60%"
$ws.Range("E7").Value2 = "2026-02-15 20:18:48"
$ws.Range("E8").Value2 = "2026-02-15 20:18:50"
$ws.Range("O8").Value2 = "8.2 °C"
$ws.Range("E9").Value2 = "2026-02-15 20:18:53"
$ws.Range("H9").Value2 = "'51%"
$ws.Range("O9").Value2 = "10.9 °C"
$ws.Range("E10").Value2 = "2026-02-15 20:18:55"
$ws.Range("E11").Value2 = "2026-02-15 20:18:58"
$ws.Range("H11").Value2 = "'43%"
$ws.Range("O11").Value2 = "7.3 °C"
$ws.Range("E12").Value2 = "2026-02-15 20:19:00"
$ws.Range("H12").Value2 = "'56%"
$ws.Range("E13").Value2 = "2026-02-15 20:19:03"
$ws.Range("H13").Value2 = "'36%"
$ws.Range("O13").Value2 = "6.4 °C"
$ws.Range("E14").Value2 = "2026-02-15 20:19:05"
$ws.Range("H14").Value2 = "'60%"
$ws.Range("E15").Value2 = "2026-02-15 20:19:06"
$ws.Range("H15").Value2 = "'51%"
$ws.Range("O15").Value2 = "10.6 °C"
$ws.Range("E16").Value2 = "2026-02-15 20:19:07"
$ws.Range("E17").Value2 = "2026-02-15 20:19:08"
$ws.Range("H17").Value2 = "'38%"
$ws.Range("E18").Value2 = "2026-02-15 20:19:10"
$ws.Range("H18").Value2 = "'72%"
$ws.Range("O18").Value2 = "7.5 °C"
$ws.Range("E19").Value2 = "2026-02-15 20:19:11"
$ws.Range("O19").Value2 = "3.4 °C"
$ws.Range("E20").Value2 = "2026-02-15 20:19:12"
$ws.Range("O20").Value2 = "-2.7 °C"
$ws.Range("E21").Value2 = "2026-02-15 20:19:13"
$ws.Range("J21").Value2 = "1015.0 hPa"
$ws.Range("E22").Value2 = "2026-02-15 20:19:14"
$ws.Range("E23").Value2 = "2026-02-15 20:19:16"
$ws.Range("I23").Value2 = "2.8 mm"
$ws.Range("O23").Value2 = "-3.6 °C"
$ws.Range("E24").Value2 = "2026-02-15 20:19:19"
$ws.Range("O24").Value2 = "8.9 °C"
$ws.Range("E25").Value2 = "2026-02-15 20:19:21"
$ws.Range("H25").Value2 = "'63%"
$ws.Range("E26").Value2 = "2026-02-15 20:19:24"
$ws.Range("E27").Value2 = "2026-02-15 20:19:26"
$ws.Range("H27").Value2 = "'49%"
$ws.Range("O27").Value2 = "0.0 °C"
$ws.Range("E28").Value2 = "2026-02-15 20:19:29"
$ws.Range("H28").Value2 = "'58%"
$ws.Range("E29").Value2 = "2026-02-15 20:19:31"
$ws.Range("H29").Value2 = "'58%"
$ws.Range("E30").Value2 = "2026-02-15 20:19:33"
$ws.Range("H30").Value2 = "'55%"
$ws.Range("E31").Value2 = "2026-02-15 20:19:35"
$ws.Range("O31").Value2 = "10.0 °C"
$ws.Range("E32").Value2 = "2026-02-15 20:19:38"
$ws.Range("L32").Value2 = "44.3 km/h - 284º 19:47 TU"
$ws.Range("E33").Value2 = "2026-02-15 20:19:40"
$ws.Range("O33").Value2 = "5.8 °C"
$ws.Range("E34").Value2 = "2026-02-15 20:19:43"
$ws.Range("H34").Value2 = "'51%"
$ws.Range("O34").Value2 = "1.1 °C"
$ws.Range("E35").Value2 = "2026-02-15 20:19:45"
$ws.Range("H35").Value2 = "'69%"
$ws.Range("E36").Value2 = "2026-02-15 20:19:48"
$ws.Range("H36").Value2 = "'48%"
$ws.Range("K36").Value2 = "10.0 MJ/m2"
$ws.Range("E37").Value2 = "2026-02-15 20:19:50"
$ws.Range("H37").Value2 = "'54%"
$ws.Range("O37").Value2 = "6.0 °C"
$ws.Range("E38").Value2 = "2026-02-15 20:19:53"
$ws.Range("O38").Value2 = "7.8 °C"
$ws.Range("E39").Value2 = "2026-02-15 20:19:55"
$ws.Range("O39").Value2 = "-2.8 °C"
$ws.Range("E40").Value2 = "2026-02-15 20:19:58"
$ws.Range("H40").Value2 = "'37%"
$ws.Range("J40").Value2 = "1016.1 hPa"
$ws.Range("E41").Value2 = "2026-02-15 20:20:00"
$ws.Range("J41").Value2 = "1016.5 hPa"
$ws.Range("K41").Value2 = "12.4 MJ/m2"
$ws.Range("O41").Value2 = "12.5 °C"
$ws.Range("E42").Value2 = "2026-02-15 20:20:02"
$ws.Range("H42").Value2 = "'57%"
$ws.Range("E43").Value2 = "2026-02-15 20:20:05"
$ws.Range("H43").Value2 = "'68%"
$ws.Range("E44").Value2 = "2026-02-15 20:20:07"
$ws.Range("I44").Value2 = "3.2 mm"
$ws.Range("O44").Value2 = "-4.0 °C"
$ws.Range("E45").Value2 = "2026-02-15 20:20:10"
$ws.Range("H45").Value2 = "'90%"
$ws.Range("I45").Value2 = "1.8 mm"
$ws.Range("E46").Value2 = "2026-02-15 20:20:12"
